$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unnecessary "nome" (name) column (column B), shifting
# everything to its right one column to the left.
$ws.Range("B:B").Delete()

$ws.Range("F12").Select()
